$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2021-09-16", "overview", "K02000001", "United Kingdom", 7339009, 26911, 158, 134805),
    @("2021-09-17", "overview", "K02000001", "United Kingdom", 7371301, 32651, 178, 134983),
    @("2021-09-18", "overview", "K02000001", "United Kingdom", 7400739, 30144, 164, 135147),
    @("2021-09-19", "overview", "K02000001", "United Kingdom", 7429746, 29612, 56, 135203),
    @("2021-09-20", "overview", "K02000001", "United Kingdom", 7465448, 36100, 49, 135252)
)

$startRow = 401
$endRow = $startRow + $data.Count - 1

# Prevent Excel from auto-parsing the date-like text in column A as a real date;
# the source data stores it as plain text, so force the column to Text format first.
$dateRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]
}
